$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.29000000000051
$ws.Range("H2").Value = 0.0001329492889827755
$ws.Range("I2").Value = 0.0001329492889827755
$ws.Range("L2").Value = 43.32340847890948
$ws.Range("M2").Value = "[20.970947141883414, 65.67586981593554]"
$ws.Range("N2").Value = 0.000314305368861012
$ws.Range("O2").Value = 0.000314305368861012
$ws.Range("P2").Value = 1.704447665911579
$ws.Range("Q2").Value = "[1.0503422885875784, 2.35855304323558]"
$ws.Range("R2").Value = [double]"4.007484533552486e-06"
$ws.Range("S2").Value = [double]"4.007484533552486e-06"
$ws.Range("T2").Value = 50.81462721899463
$ws.Range("U2").Value = "[37.26871113755763, 64.36054330043163]"
$ws.Range("V2").Value = [double]"1.540767513574792e-09"
$ws.Range("W2").Value = [double]"1.540767513574792e-09"
$ws.Range("X2").Value = 18.42954954954993
$ws.Range("Y2").Value = 15.79675675675708
$ws.Range("Z2").Value = 21.06234234234277
$ws.Range("F3").Value = 25.29000000000051
$ws.Range("H3").Value = 0.009404625739463546
$ws.Range("I3").Value = 0.009404625739463546
$ws.Range("L3").Value = 31.18663509599834
$ws.Range("M3").Value = "[8.261100068970208, 54.11217012302648]"
$ws.Range("N3").Value = 0.008780231703852781
$ws.Range("O3").Value = 0.008780231703852781
$ws.Range("P3").Value = 1.566079220708426
$ws.Range("Q3").Value = "[0.5346053564667326, 2.5975530849501194]"
$ws.Range("R3").Value = 0.003744111474839329
$ws.Range("S3").Value = 0.003744111474839329
$ws.Range("T3").Value = 57.9684014176738
$ws.Range("U3").Value = "[43.81406802493275, 72.12273481041485]"
$ws.Range("V3").Value = [double]"1.502287183541284e-10"
$ws.Range("W3").Value = [double]"1.502287183541284e-10"
$ws.Range("X3").Value = 18.98648648648687
$ws.Range("Y3").Value = 14.83477477477507
$ws.Range("Z3").Value = 23.13819819819866
$ws.Range("F4").Value = 25.29000000000051
$ws.Range("H4").Value = 0.01459868353999727
$ws.Range("I4").Value = 0.01459868353999727
$ws.Range("L4").Value = 26.29223080446661
$ws.Range("M4").Value = "[4.5131580822661235, 48.07130352666709]"
$ws.Range("N4").Value = 0.01908235934236924
$ws.Range("O4").Value = 0.01908235934236924
$ws.Range("P4").Value = 1.402552876377425
$ws.Range("Q4").Value = "[0.24528951649649944, 2.5598162362583503]"
$ws.Range("R4").Value = 0.01864471585941518
$ws.Range("S4").Value = 0.01864471585941518
$ws.Range("T4").Value = 49.50185021631177
$ws.Range("U4").Value = "[36.972991832747894, 62.03070859987565]"
$ws.Range("V4").Value = [double]"3.971272199976283e-10"
$ws.Range("W4").Value = [double]"3.971272199976283e-10"
$ws.Range("X4").Value = 19.64468468468509
$ws.Range("Y4").Value = 14.98666666666697
$ws.Range("Z4").Value = 24.3027027027032
$ws.Range("F5").Value = 25.29000000000051
$ws.Range("H5").Value = [double]"5.741969242523126e-06"
$ws.Range("I5").Value = [double]"5.741969242523126e-06"
$ws.Range("L5").Value = 52.38558132301438
$ws.Range("M5").Value = "[26.915910940462908, 77.85525170556585]"
$ws.Range("N5").Value = 0.0001492541438408157
$ws.Range("O5").Value = 0.0001492541438408157
$ws.Range("P5").Value = 0.6100790538502698
$ws.Range("Q5").Value = "[0.14465791998511524, 1.0755001877154244]"
$ws.Range("R5").Value = 0.01135049706026203
$ws.Range("S5").Value = 0.01135049706026203
$ws.Range("T5").Value = 74.64030306970081
$ws.Range("U5").Value = "[61.69832880257721, 87.5822773368244]"
$ws.Range("V5").Value = [double]"3.774758283725532e-15"
$ws.Range("W5").Value = [double]"3.774758283725532e-15"
$ws.Range("X5").Value = 22.83441441441488
$ws.Range("Y5").Value = 20.96108108108151
$ws.Range("Z5").Value = 24.70774774774825
$ws.Range("F6").Value = 24.50000000000039
$ws.Range("H6").Value = 0.001038573859302283
$ws.Range("I6").Value = 0.001038573859302283
$ws.Range("L6").Value = 43.09807545013506
$ws.Range("M6").Value = "[15.179590712389128, 71.01656018788098]"
$ws.Range("N6").Value = 0.003249514089439298
$ws.Range("O6").Value = 0.003249514089439298
$ws.Range("P6").Value = 0.3459211130078854
$ws.Range("Q6").Value = "[-0.3710790121357306, 1.0629212381515014]"
$ws.Range("R6").Value = 0.3363859907714357
$ws.Range("S6").Value = 0.3363859907714357
$ws.Range("T6").Value = 59.16735452444348
$ws.Range("U6").Value = "[44.01726378697283, 74.31744526191413]"
$ws.Range("V6").Value = [double]"5.406519676398602e-10"
$ws.Range("W6").Value = [double]"5.406519676398602e-10"
$ws.Range("X6").Value = 23.15115115115152
$ws.Range("Y6").Value = 20.35535535535568
$ws.Range("Z6").Value = 25.94694694694736
$ws.Range("F7").Value = 24.50000000000039
$ws.Range("H7").Value = [double]"1.143063180442105e-06"
$ws.Range("I7").Value = [double]"1.143063180442105e-06"
$ws.Range("L7").Value = 49.60846321346344
$ws.Range("M7").Value = "[28.864415454180275, 70.3525109727466]"
$ws.Range("N7").Value = [double]"1.689566900431316e-05"
$ws.Range("O7").Value = [double]"1.689566900431316e-05"
$ws.Range("P7").Value = 0.3710790121357315
$ws.Range("Q7").Value = "[-0.09434212172942313, 0.8365001460008861]"
$ws.Range("R7").Value = 0.115304902713121
$ws.Range("S7").Value = 0.115304902713121
$ws.Range("T7").Value = 63.81519396921546
$ws.Range("U7").Value = "[52.42263138580792, 75.20775655262301]"
$ws.Range("V7").Value = [double]"1.043609643147647e-14"
$ws.Range("W7").Value = [double]"1.043609643147647e-14"
$ws.Range("X7").Value = 23.05305305305342
$ws.Range("Y7").Value = 21.23823823823858
$ws.Range("Z7").Value = 24.86786786786827
$ws.Range("F8").Value = 24.50000000000039
$ws.Range("H8").Value = 0.07216644972170105
$ws.Range("I8").Value = 0.07216644972170105
$ws.Range("L8").Value = 25.30883776526227
$ws.Range("M8").Value = "[-2.81639876527678, 53.43407429580132]"
$ws.Range("N8").Value = 0.07659699795691832
$ws.Range("O8").Value = 0.07659699795691832
$ws.Range("P8").Value = 1.33965812855781
$ws.Range("Q8").Value = "[-0.20755266780473036, 2.88686892492035]"
$ws.Range("R8").Value = 0.08800162604420558
$ws.Range("S8").Value = 0.08800162604420558
$ws.Range("T8").Value = 61.60309217120391
$ws.Range("U8").Value = "[46.21946266679163, 76.9867216756162]"
$ws.Range("V8").Value = [double]"2.76928258102771e-10"
$ws.Range("W8").Value = [double]"2.76928258102771e-10"
$ws.Range("X8").Value = 19.27627627627658
$ws.Range("Y8").Value = 13.24324324324346
$ws.Range("Z8").Value = 25.30930930930971
$ws.Range("F9").Value = 24.50000000000039
$ws.Range("H9").Value = [double]"3.816930207678482e-06"
$ws.Range("I9").Value = [double]"3.816930207678482e-06"
$ws.Range("L9").Value = 53.22940685596252
$ws.Range("M9").Value = "[32.32779811790684, 74.13101559401821]"
$ws.Range("N9").Value = [double]"5.974861014168553e-06"
$ws.Range("O9").Value = [double]"5.974861014168553e-06"
$ws.Range("P9").Value = -0.1257894956392311
$ws.Range("Q9").Value = "[-0.6163685286322318, 0.3647895373537695]"
$ws.Range("R9").Value = 0.6080781158830737
$ws.Range("S9").Value = 0.6080781158830737
$ws.Range("T9").Value = 61.36377874173499
$ws.Range("U9").Value = "[48.36113933582753, 74.36641814764245]"
$ws.Range("V9").Value = [double]"2.487121619765276e-12"
$ws.Range("W9").Value = [double]"2.487121619765276e-12"
$ws.Range("X9").Value = 0.490490490490501
$ws.Range("Y9").Value = -1.422422422422442
$ws.Range("Z9").Value = 2.403403403403444
